$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the B:E columns (16-rep / 20-rep Lichtwark values) with the
# corresponding Fuku 16/20 columns (O, R, AN, AQ), overwriting the old
# deleted values per "Hjemme passive tweaks lichtwark deleted values".
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 44.300831108300457
$ws.Range("C2").Value = 39.221224375141411
$ws.Range("D2").Value = 48.662927375082532
$ws.Range("E2").Value = 28.480521867032515

$ws.Range("B3").Value = 68.091621566852538
$ws.Range("C3").Value = 37.859936733572958
$ws.Range("D3").Value = 45.27787645089181
$ws.Range("E3").Value = 17.934304638059899

# Update the selection to match the new, smaller block of interest.
$ws.Range("B1:E3").Select()
